$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.580.62"
$ws.Range("E2").Value = "  +4.51%  "
$ws.Range("D3").Value = "2.583.50"
$ws.Range("E3").Value = "  +5.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'587.60"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").Value = "'155.41"
$ws.Range("E6").Value = "  +6.30%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "'0.545"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").Value = "2.605.90"
$ws.Range("E9").Value = "  +6.18%  "
$ws.Range("D10").Value = "'0.114"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.32"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.361"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "'29.14"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "3.059.58"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").Value = "65.622.10"
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").Value = "2.613.20"
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("E19").Value = "  +7.43%  "
$ws.Range("D20").Value = "'11.19"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("D21").Value = "'355.35"
$ws.Range("E21").Value = "  +10.57%  "
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").Value = "'0.997"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'10.05"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").Value = "'66.18"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").Value = "'636.09"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = "  +9.76%  "
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +8.62%  "
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  +7.26%  "
$ws.Range("E38").Value = "  +6.44%  "
$ws.Range("D39").Value = "'19.35"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").Value = "'2.87"
$ws.Range("E40").Value = "  +5.76%  "
$ws.Range("D41").Value = "'155.05"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("E43").Value = "  +6.03%  "
$ws.Range("E44").Value = "  +6.69%  "
$ws.Range("D45").Value = "'42.02"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'162.49"
$ws.Range("E46").Value = "  +6.44%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'16.15"
$ws.Range("E48").Value = "  +4.81%  "
$ws.Range("D49").Value = "'3.77"
$ws.Range("E49").Value = "  +7.04%  "
$ws.Range("D50").Value = "'21.71"
$ws.Range("E50").Value = "  +7.98%  "
$ws.Range("E51").Value = "  +5.06%  "
